$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap match data (columns F:V) between paired rows while keeping
# Indice/pais/torneio/temporada/data_partida (A:E) on their original row.
$swapPairs = @(
    @(16, 17),
    @(18, 19),
    @(24, 25),
    @(31, 32)
)
foreach ($pair in $swapPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    $range1 = $ws.Range("F$r1" + ":V$r1")
    $range2 = $ws.Range("F$r2" + ":V$r2")
    $v1 = $range1.Value()
    $v2 = $range2.Value()
    $range1.Value = $v2
    $range2.Value = $v1
}

# Append 5 newly scraped fixtures (rows 46-50) at the bottom of the table.
# Row 46
$ws.Range("A45").Copy($ws.Range("A46"))
$ws.Range("E45").Copy($ws.Range("E46"))
$ws.Range("A46").Value = 45
$ws.Range("B46").Value = "montenegro"
$ws.Range("C46").Value = "prva-crnogorska-liga"
$ws.Range("D46").Value = "2023-2024"
$ws.Range("E46").Value = 45196.64583333334
$ws.Range("F46").Value = "Arsenal Tivat"
$ws.Range("G46").Value = 2
$ws.Range("H46").Value = "Rudar"
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 1.98
$ws.Range("K46").Value = "26/09/2023 02:42"
$ws.Range("L46").Value = 2.01
$ws.Range("M46").Value = "27/09/2023 15:21"
$ws.Range("N46").Value = 2.96
$ws.Range("O46").Value = "26/09/2023 02:42"
$ws.Range("P46").Value = 3.18
$ws.Range("Q46").Value = "27/09/2023 15:21"
$ws.Range("R46").Value = 3.68
$ws.Range("S46").Value = "26/09/2023 02:42"
$ws.Range("T46").Value = 3.87
$ws.Range("U46").Value = "27/09/2023 15:21"
$ws.Range("V46").Value = "https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/arsenal-tivat-rudar/0xIBXS11/"

# Row 47
$ws.Range("A45").Copy($ws.Range("A47"))
$ws.Range("E45").Copy($ws.Range("E47"))
$ws.Range("A47").Value = 46
$ws.Range("B47").Value = "montenegro"
$ws.Range("C47").Value = "prva-crnogorska-liga"
$ws.Range("D47").Value = "2023-2024"
$ws.Range("E47").Value = 45196.75
$ws.Range("F47").Value = "Jedinstvo"
$ws.Range("G47").Value = 0
$ws.Range("H47").Value = "Decic"
$ws.Range("I47").Value = 1
$ws.Range("J47").Value = 3.17
$ws.Range("K47").Value = "26/09/2023 05:12"
$ws.Range("L47").Value = 5.37
$ws.Range("M47").Value = "27/09/2023 17:40"
$ws.Range("N47").Value = 2.97
$ws.Range("O47").Value = "26/09/2023 05:12"
$ws.Range("P47").Value = 3.52
$ws.Range("Q47").Value = "27/09/2023 17:40"
$ws.Range("R47").Value = 2.16
$ws.Range("S47").Value = "26/09/2023 05:12"
$ws.Range("T47").Value = 1.66
$ws.Range("U47").Value = "27/09/2023 17:40"
$ws.Range("V47").Value = "https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/jedinstvo-decic/6FI7Y8ne/"

# Row 48
$ws.Range("A45").Copy($ws.Range("A48"))
$ws.Range("E45").Copy($ws.Range("E48"))
$ws.Range("A48").Value = 47
$ws.Range("B48").Value = "montenegro"
$ws.Range("C48").Value = "prva-crnogorska-liga"
$ws.Range("D48").Value = "2023-2024"
$ws.Range("E48").Value = 45196.75
$ws.Range("F48").Value = "Sutjeska"
$ws.Range("G48").Value = 1
$ws.Range("H48").Value = "Jezero"
$ws.Range("I48").Value = 1
$ws.Range("J48").Value = 1.46
$ws.Range("K48").Value = "26/09/2023 05:12"
$ws.Range("L48").Value = 1.54
$ws.Range("M48").Value = "27/09/2023 17:31"
$ws.Range("N48").Value = 3.81
$ws.Range("O48").Value = "26/09/2023 05:12"
$ws.Range("P48").Value = 3.92
$ws.Range("Q48").Value = "27/09/2023 17:31"
$ws.Range("R48").Value = 5.94
$ws.Range("S48").Value = "26/09/2023 05:12"
$ws.Range("T48").Value = 6
$ws.Range("U48").Value = "27/09/2023 17:43"
$ws.Range("V48").Value = "https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/sutjeska-jezero/vTGFWnH7/"

# Row 49
$ws.Range("A45").Copy($ws.Range("A49"))
$ws.Range("E45").Copy($ws.Range("E49"))
$ws.Range("A49").Value = 48
$ws.Range("B49").Value = "montenegro"
$ws.Range("C49").Value = "prva-crnogorska-liga"
$ws.Range("D49").Value = "2023-2024"
$ws.Range("E49").Value = 45196.77083333334
$ws.Range("F49").Value = "Buducnost"
$ws.Range("G49").Value = 2
$ws.Range("H49").Value = "Mladost DG"
$ws.Range("I49").Value = 2
$ws.Range("J49").Value = 1.37
$ws.Range("K49").Value = "26/09/2023 05:42"
$ws.Range("L49").Value = 1.23
$ws.Range("M49").Value = "27/09/2023 18:09"
$ws.Range("N49").Value = 4.22
$ws.Range("O49").Value = "26/09/2023 05:42"
$ws.Range("P49").Value = 5.94
$ws.Range("Q49").Value = "27/09/2023 18:09"
$ws.Range("R49").Value = 6.62
$ws.Range("S49").Value = "26/09/2023 05:42"
$ws.Range("T49").Value = 11.09
$ws.Range("U49").Value = "27/09/2023 18:09"
$ws.Range("V49").Value = "https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/buducnost-mladost-dg/fuF3ZlXl/"

# Row 50
$ws.Range("A45").Copy($ws.Range("A50"))
$ws.Range("E45").Copy($ws.Range("E50"))
$ws.Range("A50").Value = 49
$ws.Range("B50").Value = "montenegro"
$ws.Range("C50").Value = "prva-crnogorska-liga"
$ws.Range("D50").Value = "2023-2024"
$ws.Range("E50").Value = 45196.79166666666
$ws.Range("F50").Value = "Mornar Bar"
$ws.Range("G50").Value = 0
$ws.Range("H50").Value = "Petrovac"
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 2.4
$ws.Range("K50").Value = "26/09/2023 12:42"
$ws.Range("L50").Value = 2.66
$ws.Range("M50").Value = "27/09/2023 18:38"
$ws.Range("N50").Value = 2.81
$ws.Range("O50").Value = "26/09/2023 12:42"
$ws.Range("P50").Value = 2.56
$ws.Range("Q50").Value = "27/09/2023 18:38"
$ws.Range("R50").Value = 2.93
$ws.Range("S50").Value = "26/09/2023 12:42"
$ws.Range("T50").Value = 3.28
$ws.Range("U50").Value = "27/09/2023 18:38"
$ws.Range("V50").Value = "https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/mornar-bar-petrovac/CGPbzVHr/"
